# PF ini generation number auto + matplotlib version
#
# 1. Duplicate the "map (2)" sheet, inserting the copy right before it,
#    then rename the copy to "map (1)".
# 2. On the new "map (1)" sheet, replace several per-component labels
#    (C0/C1/C2/C3, P0/P1/P2, E0/E1/E2/E4) with the generic "C"/"P"/"E"
#    labels, and move the selection to C5.
# 3. On the "map (5)" sheet, clear the stray "P3" label in P5 back to a
#    plain 0, move the selection to R9, and make it the active sheet/tab.

$wb = $excel.ActiveWorkbook

$src = $wb.Worksheets.Item("map (2)")
$src.Copy($src)
$ws = $wb.Worksheets.Item("map (2) (2)")
$ws.Name = "map (1)"

$ws.Range("D1").Value = "C"
$ws.Range("A3").Value = "C"
$ws.Range("C3").Value = "P"
$ws.Range("E3").Value = "E"
$ws.Range("G3").Value = "C"
$ws.Range("B5").Value = "E"
$ws.Range("C5").Value = "P"
$ws.Range("F5").Value = "E"
$ws.Range("E6").Value = "E"
$ws.Range("B8").Value = "E"
$ws.Range("C8").Value = "P"
$ws.Range("D8").Value = "E"
$ws.Range("E8").Value = "P"
$ws.Range("D10").Value = "C"

[void]$ws.Range("C5").Select()

$map5 = $wb.Worksheets.Item("map (5)")
$map5.Activate()
$map5.Range("P5").Value = 0
[void]$map5.Range("R9").Select()
